$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarterly row: 01-07-2021
# Write the date label first. Excel would normally auto-detect the
# dd-mm-yyyy text as a real date and reformat the cell; force it to stay
# plain text (matching the rest of column A) and then drop the
# temporary formatting so no new cell style is introduced.
$ws.Cells.Item(76, 1).NumberFormat = "@"
$ws.Cells.Item(76, 1).Value = "01-07-2021"
$ws.Cells.Item(76, 1).ClearFormats()

# Fill in the new data row with numeric values (columns B..AA)
$values = @(9225,9236,0,-11,3398,2377,-1,0,2377,1476,174,1162,-1984,2124,-428,-428,-27,-18,-450,366,75,9035,1935,5584,-205,1720)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(76, 2 + $i).Value = $values[$i]
}
